$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the target paragraph: "Add body from SharePoint Cerate Item action
# output" (a typo for "Create"). We search by scanning paragraphs for the
# unique substring, which is more robust than a hard-coded paragraph index.
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -like "*Cerate Item action output*") {
        $targetIndex = $i
        break
    }
}

$p = $d.Paragraphs.Item($targetIndex)
$r = $p.Range
$pStart = $r.Start
$pEnd = $r.End

# Replace the paragraph's text (excluding the trailing paragraph mark) with
# the corrected + extended sentence.
$newText = "Add body from SharePoint Create Item action output to handle the error code and message "
$textRange = $d.Range($pStart, $pEnd - 1)
$textRange.Text = $newText

# Apply bold formatting to "SharePoint Create Item", split across three
# separate runs (matching how the source document records it): "SharePoint
# Cr" + "e" + "ate Item".
$off = $pStart + "Add body from ".Length

$b1s = $off
$off = $off + "SharePoint Cr".Length
$b1e = $off
$d.Range($b1s, $b1e).Font.Bold = 1

$b2s = $off
$off = $off + "e".Length
$b2e = $off
$d.Range($b2s, $b2e).Font.Bold = 1

$b3s = $off
$off = $off + "ate Item".Length
$b3e = $off
$d.Range($b3s, $b3e).Font.Bold = 1

# ---------------------------------------------------------------------------
# Insert the new "Example:" paragraph right after the target paragraph.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item($targetIndex)
$p.Range.InsertParagraphAfter()
$exampleIndex = $targetIndex + 1
$pExample = $d.Paragraphs.Item($exampleIndex)
$rExample = $pExample.Range
$exStart = $rExample.Start
$rExample.Text = "Example:"
$d.Range($exStart, $exStart + "Example".Length).Font.Bold = 1

# ---------------------------------------------------------------------------
# Insert the 4 code-block paragraphs ( {, status line, message line, } )
# after "Example:". Each one needs the style reset to "Normal" so the
# inherited 720-twip indent is removed, then the Consolas code-block
# shading/spacing is applied.
# ---------------------------------------------------------------------------
function Format-CodeParagraph($para) {
    $para.Style = "Normal"
    $para.Format.Shading.Texture = 0
    $para.Format.Shading.ForegroundPatternColor = -16777216
    $para.Format.Shading.BackgroundPatternColor = 16711679
    $para.Format.SpaceAfter = 0
    $para.Format.LineSpacingRule = 3
    $para.Format.LineSpacing = 14.25
}

# --- "{" paragraph ---
$pExample = $d.Paragraphs.Item($exampleIndex)
$pExample.Range.InsertParagraphAfter()
$braceOpenIndex = $exampleIndex + 1
$pBraceOpen = $d.Paragraphs.Item($braceOpenIndex)
Format-CodeParagraph $pBraceOpen
$rBraceOpen = $pBraceOpen.Range
$rBraceOpen.Text = "{"
$rBraceOpen.Font.Name = "Consolas"
$rBraceOpen.Font.NameAscii = "Consolas"
$rBraceOpen.Font.Color = 0
$rBraceOpen.Font.Size = 10.5

# --- status line paragraph ---
$pBraceOpen = $d.Paragraphs.Item($braceOpenIndex)
$pBraceOpen.Range.InsertParagraphAfter()
$statusIndex = $braceOpenIndex + 1
$pStatus = $d.Paragraphs.Item($statusIndex)
Format-CodeParagraph $pStatus
$rStatus = $pStatus.Range
$statusStart = $rStatus.Start
$rStatus.Text = "  `"status`": 404,"
$rStatus2 = $pStatus.Range
$rStatus2.Font.Name = "Consolas"
$rStatus2.Font.NameAscii = "Consolas"
$rStatus2.Font.Color = 0
$rStatus2.Font.Size = 10.5

$off = $statusStart + "  ".Length
$s1 = $off
$off = $off + '"status"'.Length
$s1e = $off
$d.Range($s1, $s1e).Font.Color = 1381795

$off2 = $s1e + ": ".Length
$numStart = $off2
$numEnd = $numStart + "404".Length
$d.Range($numStart, $numEnd).Font.Color = 5933065

# --- message line paragraph ---
$pStatus = $d.Paragraphs.Item($statusIndex)
$pStatus.Range.InsertParagraphAfter()
$messageIndex = $statusIndex + 1
$pMessage = $d.Paragraphs.Item($messageIndex)
Format-CodeParagraph $pMessage
$rMessage = $pMessage.Range
$messageStart = $rMessage.Start
$msgText = '"List not found\r\nclientRequestId: 01cd9b0e-74bf-4383-880d-2f313cdd72db\r\nserviceRequestId: 01cd9b0e-74bf-4383-880d-2f313cdd72db"'
$rMessage.Text = "  `"message`": " + $msgText
$rMessage2 = $pMessage.Range
$rMessage2.Font.Name = "Consolas"
$rMessage2.Font.NameAscii = "Consolas"
$rMessage2.Font.Color = 0
$rMessage2.Font.Size = 10.5

$off = $messageStart + "  ".Length
$m1 = $off
$off = $off + '"message"'.Length
$m1e = $off
$d.Range($m1, $m1e).Font.Color = 1381795

$off2 = $m1e + ": ".Length
$valStart = $off2
$valEnd = $valStart + $msgText.Length
$d.Range($valStart, $valEnd).Font.Color = 10834180

# --- "}" paragraph ---
$pMessage = $d.Paragraphs.Item($messageIndex)
$pMessage.Range.InsertParagraphAfter()
$braceCloseIndex = $messageIndex + 1
$pBraceClose = $d.Paragraphs.Item($braceCloseIndex)
Format-CodeParagraph $pBraceClose
$rBraceClose = $pBraceClose.Range
$rBraceClose.Text = "}"
$rBraceClose.Font.Name = "Consolas"
$rBraceClose.Font.NameAscii = "Consolas"
$rBraceClose.Font.Color = 0
$rBraceClose.Font.Size = 10.5

# ---------------------------------------------------------------------------
# Final blank paragraph after the closing "}" (matching the style used by
# "Example:" / the body paragraph -- 720-twip indent, Times New Roman look).
# ---------------------------------------------------------------------------
$pBraceClose = $d.Paragraphs.Item($braceCloseIndex)
$pBraceClose.Range.InsertParagraphAfter()

Write-Host "Done"
